$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without Excel
# auto-converting numeric-looking / percent-looking strings into
# Number/Date types, and without leaving a stray number-format
# style behind on the cell (matches original inlineStr text cells).
function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# --- Rows where only price (D) changed ---
Set-TextCell 9 4 "0.999"

# --- Rows where only volume/percent (E) changed ---
Set-TextCell 20 5 "  -0.28%  "

# --- Rows where both price (D) and percent (E) changed ---
Set-TextCell 2 4 "69.008.53"
Set-TextCell 2 5 "  +4.96%  "
Set-TextCell 3 4 "3.535.06"
Set-TextCell 3 5 "  +3.71%  "
Set-TextCell 4 4 "0.998"
Set-TextCell 4 5 "  -0.06%  "
Set-TextCell 5 4 "589.08"
Set-TextCell 5 5 "  +4.94%  "
Set-TextCell 6 4 "193.43"
Set-TextCell 6 5 "  +9.46%  "
Set-TextCell 7 4 "0.638"
Set-TextCell 7 5 "  +0.70%  "
Set-TextCell 8 4 "3.527.97"
Set-TextCell 8 5 "  +3.79%  "
Set-TextCell 10 4 "0.179"
Set-TextCell 10 5 "  +3.79%  "
Set-TextCell 11 4 "0.657"
Set-TextCell 11 5 "  +2.03%  "
Set-TextCell 12 4 "58.87"
Set-TextCell 12 5 "  +9.59%  "
Set-TextCell 13 4 "0.0000292"
Set-TextCell 13 5 "  +4.71%  "
Set-TextCell 14 4 "9.62"
Set-TextCell 14 5 "  +3.93%  "
Set-TextCell 15 4 "4.059.69"
Set-TextCell 15 5 "  +3.11%  "
Set-TextCell 16 4 "19.20"
Set-TextCell 16 5 "  +4.59%  "
Set-TextCell 17 4 "3.522.14"
Set-TextCell 17 5 "  +4.11%  "
Set-TextCell 18 4 "68.682.93"
Set-TextCell 18 5 "  +4.80%  "
Set-TextCell 19 4 "12.36"
Set-TextCell 19 5 "  +4.00%  "
Set-TextCell 21 4 "1.04"
Set-TextCell 21 5 "  +2.54%  "
Set-TextCell 22 4 "491.47"
Set-TextCell 22 5 "  +0.00%  "
Set-TextCell 23 4 "5.67"
Set-TextCell 23 5 "  +14.35%  "
Set-TextCell 24 4 "17.22"
Set-TextCell 24 5 "  +21.50%  "
Set-TextCell 25 4 "4.48"
Set-TextCell 25 5 "  +8.39%  "
Set-TextCell 26 4 "90.75"
Set-TextCell 26 5 "  +1.71%  "
Set-TextCell 27 4 "3.04"
Set-TextCell 27 5 "  +3.97%  "
Set-TextCell 28 4 "11.18"
Set-TextCell 28 5 "  +4.02%  "
Set-TextCell 29 4 "9.22"
Set-TextCell 29 5 "  +4.98%  "
Set-TextCell 30 4 "31.90"
Set-TextCell 30 5 "  +1.07%  "
Set-TextCell 31 4 "7.49"
Set-TextCell 31 5 "  +13.96%  "
Set-TextCell 32 4 "614.83"
Set-TextCell 32 5 "  +6.46%  "
Set-TextCell 33 4 "11.97"
Set-TextCell 33 5 "  +3.78%  "
Set-TextCell 34 4 "65.22"
Set-TextCell 34 5 "  +3.87%  "
Set-TextCell 35 4 "0.114"
Set-TextCell 35 5 "  +4.72%  "
Set-TextCell 38 4 "37.55"
Set-TextCell 38 5 "  +4.04%  "
Set-TextCell 39 4 "0.395"
Set-TextCell 39 5 "  +5.34%  "
Set-TextCell 40 4 "0.0₃0789"
Set-TextCell 40 5 "  +5.99%  "
Set-TextCell 41 4 "3.53"
Set-TextCell 41 5 "  -2.97%  "
Set-TextCell 42 4 "3.268.80"
Set-TextCell 42 5 "  +4.64%  "
Set-TextCell 43 4 "2.96"
Set-TextCell 43 5 "  +5.87%  "
Set-TextCell 44 4 "0.0439"
Set-TextCell 44 5 "  +4.74%  "
Set-TextCell 45 4 "2.63"
Set-TextCell 45 5 "  +7.49%  "
Set-TextCell 46 4 "3.31"
Set-TextCell 46 5 "  +3.82%  "
Set-TextCell 49 4 "9.02"
Set-TextCell 49 5 "  +6.34%  "
Set-TextCell 50 4 "0.998"
Set-TextCell 50 5 "  +0.01%  "
Set-TextCell 51 4 "141.02"
Set-TextCell 51 5 "  +0.38%  "

# --- Rows 36/37 and 47/48 swapped (coin name, link, price, percent) ---
Set-TextCell 36 2 "Kaspa"
Set-TextCell 36 3 "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell 36 4 "0.148"
Set-TextCell 36 5 "  +4.76%  "
Set-TextCell 37 2 "Dai"
Set-TextCell 37 3 "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell 37 4 "1.00"
Set-TextCell 37 5 "  -0.03%  "
Set-TextCell 47 2 "dogwifhat"
Set-TextCell 47 3 "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell 47 4 "2.80"
Set-TextCell 47 5 "  +20.19%  "
Set-TextCell 48 2 "Stellar"
Set-TextCell 48 3 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell 48 4 "0.137"
Set-TextCell 48 5 "  +1.29%  "
